$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(14, 2).Value = 0.13325667381286621
$ws.Cells.Item(14, 3).Value = 0.3525640070438385
$ws.Cells.Item(15, 2).Value = 0.59630095958709717
$ws.Cells.Item(15, 3).Value = 1.2037181854248047
$ws.Cells.Item(16, 2).Value = 2.3343579769134521
$ws.Cells.Item(16, 3).Value = 1.6094293594360352
$ws.Cells.Item(17, 2).Value = 4.9510679244995117
$ws.Cells.Item(17, 3).Value = 1.3200004100799561
$ws.Cells.Item(18, 2).Value = 7.2638721466064453
$ws.Cells.Item(18, 3).Value = 1.2844691276550293
$ws.Cells.Item(19, 2).Value = 8.2130584716796875
$ws.Cells.Item(19, 3).Value = 2.0285651683807373
$ws.Cells.Item(20, 2).Value = 7.4087648391723633
$ws.Cells.Item(20, 3).Value = 2.5005242824554443
$ws.Cells.Item(21, 2).Value = 5.2815642356872559
$ws.Cells.Item(21, 3).Value = 2.2695682048797607
$ws.Cells.Item(22, 2).Value = 2.8288061618804932
$ws.Cells.Item(22, 3).Value = 1.6792465448379517
$ws.Cells.Item(23, 2).Value = 0.97551333904266357
$ws.Cells.Item(23, 3).Value = 1.3961265087127686
$ws.Cells.Item(24, 2).Value = 0.32703468203544617
$ws.Cells.Item(24, 3).Value = 0.86525249481201172
$ws.Cells.Item(25, 2).Value = 0.13718140125274658
$ws.Cells.Item(25, 3).Value = 0.36294785141944885
$ws.Cells.Item(26, 2).Value = 0.027017392218112946
$ws.Cells.Item(26, 3).Value = 0.071481294929981232
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(29, 2).Value = 0
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(30, 2).Value = 0
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(34, 2).Value = 0
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(35, 2).Value = 0
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(36, 2).Value = 0
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(37, 2).Value = 0
$ws.Cells.Item(37, 3).Value = 0
$ws.Cells.Item(38, 2).Value = 0
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(39, 2).Value = 0
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(40, 2).Value = 0
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(43, 2).Value = 0.0021573405247181654
$ws.Cells.Item(43, 3).Value = 0.0057077864184975624
$ws.Cells.Item(44, 2).Value = 0.00935140810906887
$ws.Cells.Item(44, 3).Value = 0.024741498753428459
$ws.Cells.Item(45, 2).Value = 0.015582025982439518
$ws.Cells.Item(45, 3).Value = 0.041226167231798172
$ws.Cells.Item(46, 2).Value = 0.016120001673698425
$ws.Cells.Item(46, 3).Value = 0.042649514973163605
$ws.Cells.Item(47, 2).Value = 0.011006052605807781
$ws.Cells.Item(47, 3).Value = 0.029119279235601425
$ws.Cells.Item(48, 2).Value = 0.0041294349357485771
$ws.Cells.Item(48, 3).Value = 0.010925457812845707
$ws.Cells.Item(49, 2).Value = 0
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(50, 2).Value = 0
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(51, 2).Value = 0.059799175709486008
$ws.Cells.Item(51, 3).Value = 0.1552794873714447
$ws.Cells.Item(52, 2).Value = 0.18992060422897339
$ws.Cells.Item(52, 3).Value = 0.47109448909759521
$ws.Cells.Item(53, 2).Value = 0.32295951247215271
$ws.Cells.Item(53, 3).Value = 0.76917624473571777
$ws.Cells.Item(54, 2).Value = 0.38927364349365234
$ws.Cells.Item(54, 3).Value = 0.88162171840667725
$ws.Cells.Item(55, 2).Value = 0.42728027701377869
$ws.Cells.Item(55, 3).Value = 0.72623181343078613
$ws.Cells.Item(56, 2).Value = 0.46789297461509705
$ws.Cells.Item(56, 3).Value = 0.4724128246307373
$ws.Cells.Item(57, 2).Value = 0.51833248138427734
$ws.Cells.Item(57, 3).Value = 0.47248518466949463
$ws.Cells.Item(58, 2).Value = 0.54207545518875122
$ws.Cells.Item(58, 3).Value = 0.54727208614349365
$ws.Cells.Item(59, 2).Value = 0.51892364025115967
$ws.Cells.Item(59, 3).Value = 0.45023238658905029
$ws.Cells.Item(60, 2).Value = 0.41606965661048889
$ws.Cells.Item(60, 3).Value = 0.42585000395774841
$ws.Cells.Item(61, 2).Value = 0.30833426117897034
$ws.Cells.Item(61, 3).Value = 0.55056506395339966
$ws.Cells.Item(62, 2).Value = 0.24714547395706177
$ws.Cells.Item(62, 3).Value = 0.61282801628112793
$ws.Cells.Item(63, 2).Value = 0.1987127810716629
$ws.Cells.Item(63, 3).Value = 0.52574455738067627
$ws.Cells.Item(64, 2).Value = 0.12490011006593704
$ws.Cells.Item(64, 3).Value = 0.33045461773872375
$ws.Cells.Item(65, 2).Value = 0.046370331197977066
$ws.Cells.Item(65, 3).Value = 0.12268436700105667
$ws.Cells.Item(66, 2).Value = 0.00022212395560927689
$ws.Cells.Item(66, 3).Value = 0.00058768480084836483
$ws.Cells.Item(67, 2).Value = 0.010367498733103275
$ws.Cells.Item(67, 3).Value = 0.027429824694991112
$ws.Cells.Item(68, 2).Value = 0.057072855532169342
$ws.Cells.Item(68, 3).Value = 0.098446264863014221
$ws.Cells.Item(69, 2).Value = 0.18138702213764191
$ws.Cells.Item(69, 3).Value = 0.33669793605804443
$ws.Cells.Item(70, 2).Value = 0.34199383854866028
$ws.Cells.Item(70, 3).Value = 0.69254380464553833
$ws.Cells.Item(71, 2).Value = 0.46732673048973083
$ws.Cells.Item(71, 3).Value = 0.99592804908752441
$ws.Cells.Item(72, 2).Value = 0.49863070249557495
$ws.Cells.Item(72, 3).Value = 1.1008104085922241
